$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(1.02, 1.048154189308476, 1.054302183303247, 1.059078281672559, 1.066550078969349, 1.047422977314475, 1.053199449270076, 1.057045768333732, 1.061808767481459, 1.069260326811885, 1.021545651555174)
$row3 = @(1.02, 1.04913576325401, 1.055084001901032, 1.060055359126811, 1.067524359233582, 1.047718898713166, 1.053829628130604, 1.057640883458389, 1.062599591252257, 1.07004982529948, 1.021758240382928)
$row4 = @(1.02, 1.049771058748669, 1.05559004213255, 1.060688655791143, 1.068155557672584, 1.047909133268098, 1.054236923959538, 1.058025454619695, 1.06311174795419, 1.070560834666068, 1.021895563769327)
$row5 = @(1.02, 1.050038172689308, 1.055802816652312, 1.060955146961633, 1.068421097756509, 1.047988809057039, 1.054408037454384, 1.058187005967689, 1.063327163292843, 1.070775698291258, 1.021953237771691)
$row6 = @(1.02, 1.050083024373479, 1.055838544465308, 1.060999906822353, 1.068465693880319, 1.048002169437223, 1.054436761478015, 1.058214123953135, 1.063363338588008, 1.070811776883256, 1.021962918154)
$row7 = @(1.02, 1.049774627798572, 1.055592885099039, 1.060692215665398, 1.068159105109448, 1.047910199074577, 1.054239210833914, 1.058027613758442, 1.063114625934336, 1.070563705546781, 1.021896334635528)
$row8 = @(1.02, 1.048485885110075, 1.05456637061403, 1.059408269005808, 1.066879180674625, 1.047523243146972, 1.053412518503299, 1.057246995050558, 1.062075938270435, 1.069527109719327, 1.021617545589663)
$row9 = @(1.02, 1.04621614071675, 1.052758724734131, 1.057153972289869, 1.064629763148433, 1.046831850623533, 1.051952195042402, 1.055867585380746, 1.060249052882164, 1.067701690833803, 1.021124490833441)
$row10 = @(1.02, 1.044703805285126, 1.051554496107365, 1.055656669318426, 1.063134226880043, 1.046364541861962, 1.050976273210976, 1.054945425559868, 1.059033476341922, 1.066485599292979, 1.020794599532389)
$row11 = @(1.02, 1.044049149428152, 1.051033269943737, 1.055009652147845, 1.062487621913804, 1.04616068361112, 1.050553133080819, 1.054545522958533, 1.058507685360044, 1.065959231144457, 1.020651474655175)
$row12 = @(1.02, 1.043806010826233, 1.050839696233242, 1.054769521172981, 1.062247591012464, 1.046084734997087, 1.05039587648478, 1.054396891608779, 1.058312468358256, 1.065763746528781, 1.020598269934216)
$row13 = @(1.02, 1.043858163530117, 1.050881216947835, 1.054821021015892, 1.062299071789887, 1.046101036480346, 1.050429612328627, 1.054428777587761, 1.058354339200498, 1.065805677184116, 1.020609684413379)
$row14 = @(1.02, 1.044029050906732, 1.05101726840405, 1.054989798760835, 1.062467777881055, 1.04615441029519, 1.050540135901938, 1.05453323887141, 1.058491546919541, 1.065943071668862, 1.020647077586896)
$row15 = @(1.02, 1.044134344210686, 1.051101098583019, 1.055093814831984, 1.062571742774912, 1.046187265650111, 1.050608222073689, 1.054597589046072, 1.058576096434583, 1.066027729202558, 1.02067011122831)
$row16 = @(1.02, 1.044747256786799, 1.051589092733334, 1.05569963771308, 1.063177160482212, 1.046378039468456, 1.051004343886177, 1.054971953162604, 1.0590683832594, 1.066520537051599, 1.020804092387187)
$row17 = @(1.02, 1.045131773225093, 1.051895256192379, 1.056080009790284, 1.063557183928312, 1.046497302693481, 1.051252671219579, 1.055206621332852, 1.059377332716179, 1.066829718499143, 1.020888060512412)
$row18 = @(1.02, 1.045356073982201, 1.052073856626078, 1.056302002156365, 1.063778938907036, 1.046566721152602, 1.051397462327432, 1.055343441246436, 1.05955759184812, 1.067010078776574, 1.02093701068814)
$row19 = @(1.02, 1.045432557925302, 1.052134758202033, 1.056377717459037, 1.063854567494588, 1.046590366318843, 1.051446823173209, 1.05539008343434, 1.059619064712912, 1.06707158033397, 1.020953696843532)
$row20 = @(1.02, 1.045090516300394, 1.051862405626895, 1.0560391862435, 1.063516401307169, 1.046484521952458, 1.051226033630803, 1.055181449659924, 1.059344179757527, 1.066796544184652, 1.020879054319212)
$row21 = @(1.02, 1.043978728009751, 1.050977203725703, 1.054940092392627, 1.062418094077475, 1.046138699289012, 1.050507591765201, 1.054502480098652, 1.058451140308335, 1.06590261152933, 1.02063606738021)
$row22 = @(1.02, 1.043279874795843, 1.050420832931652, 1.054250206143026, 1.061728395531547, 1.0459199560903, 1.050055395953187, 1.054075066216011, 1.057890144180024, 1.065340745833704, 1.02048305031573)
$row23 = @(1.02, 1.04365033365551, 1.050715757223828, 1.054615817841433, 1.062093936824759, 1.046036040127074, 1.050295159096963, 1.054301695367446, 1.058187491825171, 1.065638583651581, 1.020564190363375)
$row24 = @(1.02, 1.04510915847056, 1.051877249320214, 1.056057632252338, 1.063534828928625, 1.046490297469562, 1.051238070177496, 1.055192823834828, 1.059359159985243, 1.066811534167725, 1.020883123913723)
$row25 = @(1.02, 1.046802780466384, 1.053225895758421, 1.057735786234059, 1.065210577546602, 1.047011718578489, 1.052330144688975, 1.056224648951276, 1.06072093660273, 1.068173458494329, 1.021252167987374)

$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(2, $cols[$i]).Value = $row2[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(3, $cols[$i]).Value = $row3[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(4, $cols[$i]).Value = $row4[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(5, $cols[$i]).Value = $row5[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(6, $cols[$i]).Value = $row6[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(7, $cols[$i]).Value = $row7[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(8, $cols[$i]).Value = $row8[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(9, $cols[$i]).Value = $row9[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(10, $cols[$i]).Value = $row10[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(11, $cols[$i]).Value = $row11[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(12, $cols[$i]).Value = $row12[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(13, $cols[$i]).Value = $row13[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(14, $cols[$i]).Value = $row14[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(15, $cols[$i]).Value = $row15[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(16, $cols[$i]).Value = $row16[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(17, $cols[$i]).Value = $row17[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(18, $cols[$i]).Value = $row18[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(19, $cols[$i]).Value = $row19[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(20, $cols[$i]).Value = $row20[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(21, $cols[$i]).Value = $row21[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(22, $cols[$i]).Value = $row22[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(23, $cols[$i]).Value = $row23[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(24, $cols[$i]).Value = $row24[$i] }
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Cells.Item(25, $cols[$i]).Value = $row25[$i] }
